$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 70
$ws1.Range("F7").Value = 588
$ws1.Range("F8").Value = 107
$ws1.Range("F9").Value = 8703
$ws1.Range("F10").Value = 806
$ws1.Range("F12").Value = 1142
$ws1.Range("F13").Value = 977
$ws1.Range("F17").Value = 234
$ws1.Range("F18").Value = 255
$ws1.Range("F21").Value = 1018

# Sheet "全部类型" (All types) - same underlying rows, column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 70
$ws4.Range("F9").Value = 588
$ws4.Range("F10").Value = 107
$ws4.Range("F11").Value = 8703
$ws4.Range("F12").Value = 806
$ws4.Range("F14").Value = 1142
$ws4.Range("F15").Value = 977
$ws4.Range("F19").Value = 234
$ws4.Range("F20").Value = 255
$ws4.Range("F23").Value = 1018
